$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" shifts from D to E)
$ws.Columns.Item(4).Insert()

# New header "MAE" in D1, matching the style of the other header cells (row 1)
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# New MAE values in D2:D5
$ws.Range("D2").Value = 0.5256692707344185
$ws.Range("D3").Value = 0.3215706312231011
$ws.Range("D4").Value = 0.3777406617731509
$ws.Range("D5").Value = 0.468209296615762

# Minor floating point correction on B5
$ws.Range("B5").Value = 0.3638702225807679
